# Scheduled-runner data refresh for Shinryu_Profits workbook.
# Updates market-board derived price/profit columns (H:N) on a handful of
# leve rows across the ALC, BSM, CUL and WVR sheets.

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 207.5625
$ws.Range("I12").Value = 215.88889
$ws.Range("J12").Value = 196.85715
$ws.Range("K12").Value = 215.88889
$ws.Range("L12").Value = 196.85715
$ws.Range("M12").Value = -45.88889
$ws.Range("N12").Value = -536.85715
$ws.Range("H21").Value = 53346
$ws.Range("I21").Value = 53346
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 53346
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -52878
$ws.Range("H23").Value = 53346
$ws.Range("I23").Value = 53346
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 53346
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -53112
$ws.Range("H29").Value = 232.66667
$ws.Range("I29").Value = 149
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 447
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = -166
$ws.Range("N29").Value = -1762
$ws.Range("H38").Value = 1079.0667
$ws.Range("I38").Value = 199.64706
$ws.Range("J38").Value = 2229.077
$ws.Range("K38").Value = 598.94118
$ws.Range("L38").Value = 6687.231000000001
$ws.Range("M38").Value = -226.94118
$ws.Range("N38").Value = -7431.231000000001
$ws.Range("H58").Value = 1262.8
$ws.Range("I58").Value = 995.1667
$ws.Range("J58").Value = 2333.3333
$ws.Range("K58").Value = 2985.5001
$ws.Range("L58").Value = 6999.999899999999
$ws.Range("M58").Value = -2835.5001
$ws.Range("N58").Value = -7299.999899999999
$ws.Range("H87").Value = 32983.22
$ws.Range("J87").Value = 32983.22
$ws.Range("L87").Value = 32983.22
$ws.Range("N87").Value = -35479.22
$ws.Range("H90").Value = 32983.22
$ws.Range("J90").Value = 32983.22
$ws.Range("L90").Value = 98949.66
$ws.Range("N90").Value = -111429.66
$ws.Range("H103").Value = 1639
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 2100.7778
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 6302.3334
$ws.Range("M103").Value = -1214
$ws.Range("N103").Value = -7474.3334
$ws.Range("H141").Value = 1915.2
$ws.Range("I141").Value = 1457.6666
$ws.Range("K141").Value = 4372.9998
$ws.Range("M141").Value = 807.0002000000004

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2991.5151
$ws.Range("I105").Value = 2910
$ws.Range("J105").Value = 2996.7742
$ws.Range("K105").Value = 2910
$ws.Range("L105").Value = 2996.7742
$ws.Range("M105").Value = -1163
$ws.Range("N105").Value = -6490.7742

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3222.6
$ws.Range("I75").Value = 113
$ws.Range("J75").Value = 4000
$ws.Range("K75").Value = 339
$ws.Range("L75").Value = 12000
$ws.Range("M75").Value = 659
$ws.Range("N75").Value = -13996
$ws.Range("H78").Value = 3222.6
$ws.Range("I78").Value = 113
$ws.Range("J78").Value = 4000
$ws.Range("K78").Value = 1017
$ws.Range("L78").Value = 36000
$ws.Range("M78").Value = 3975
$ws.Range("N78").Value = -45984
$ws.Range("H81").Value = 45251.145
$ws.Range("I81").Value = 1338.6666
$ws.Range("J81").Value = 57227.273
$ws.Range("K81").Value = 4015.9998
$ws.Range("L81").Value = 171681.819
$ws.Range("M81").Value = -2892.9998
$ws.Range("N81").Value = -173927.819
$ws.Range("H84").Value = 45251.145
$ws.Range("I84").Value = 1338.6666
$ws.Range("J84").Value = 57227.273
$ws.Range("K84").Value = 12047.9994
$ws.Range("L84").Value = 515045.457
$ws.Range("M84").Value = -6431.999400000001
$ws.Range("N84").Value = -526277.4569999999
$ws.Range("H88").Value = 17500
$ws.Range("J88").Value = 17500
$ws.Range("L88").Value = 52500
$ws.Range("N88").Value = -53356
$ws.Range("H91").Value = 17500
$ws.Range("J91").Value = 17500
$ws.Range("L91").Value = 52500
$ws.Range("N91").Value = -55464
$ws.Range("H93").Value = 4887.4443
$ws.Range("I93").Value = 3024
$ws.Range("J93").Value = 4997.0586
$ws.Range("K93").Value = 9072
$ws.Range("L93").Value = 14991.1758
$ws.Range("M93").Value = -7200
$ws.Range("N93").Value = -18735.1758
$ws.Range("H94").Value = 3667.8333
$ws.Range("I94").Value = 1024
$ws.Range("J94").Value = 3908.182
$ws.Range("K94").Value = 3072
$ws.Range("L94").Value = 11724.546
$ws.Range("M94").Value = -2396
$ws.Range("N94").Value = -13076.546
$ws.Range("H99").Value = 1967.7273
$ws.Range("I99").Value = 1205.625
$ws.Range("K99").Value = 3616.875
$ws.Range("M99").Value = -1370.875
$ws.Range("H105").Value = 3800
$ws.Range("J105").Value = 3800
$ws.Range("L105").Value = 11400
$ws.Range("N105").Value = -16642
$ws.Range("H108").Value = 2122
$ws.Range("I108").Value = 395.2
$ws.Range("K108").Value = 1185.6
$ws.Range("M108").Value = 1694.4
$ws.Range("H110").Value = 3596.375
$ws.Range("I110").Value = 2528.5
$ws.Range("J110").Value = 6800
$ws.Range("K110").Value = 7585.5
$ws.Range("L110").Value = 20400
$ws.Range("M110").Value = -3495.5
$ws.Range("N110").Value = -28580
$ws.Range("H111").Value = 989.4
$ws.Range("I111").Value = 511.75
$ws.Range("J111").Value = 2900
$ws.Range("K111").Value = 1535.25
$ws.Range("L111").Value = 8700
$ws.Range("M111").Value = 1531.75
$ws.Range("N111").Value = -14834
$ws.Range("H112").Value = 2914.4783
$ws.Range("I112").Value = 1654.125
$ws.Range("J112").Value = 3586.6667
$ws.Range("K112").Value = 4962.375
$ws.Range("L112").Value = 10760.0001
$ws.Range("M112").Value = -3854.375
$ws.Range("N112").Value = -12976.0001
$ws.Range("H114").Value = 1611.8966
$ws.Range("I114").Value = 1112.6364
$ws.Range("J114").Value = 1917
$ws.Range("K114").Value = 3337.9092
$ws.Range("L114").Value = 5751
$ws.Range("M114").Value = -83.90920000000006
$ws.Range("N114").Value = -12259
$ws.Range("H115").Value = 2362.6086
$ws.Range("I115").Value = 1099
$ws.Range("J115").Value = 3334.6155
$ws.Range("K115").Value = 3297
$ws.Range("L115").Value = 10003.8465
$ws.Range("M115").Value = -2122
$ws.Range("N115").Value = -12353.8465
$ws.Range("H116").Value = 1973
$ws.Range("I116").Value = 1158.375
$ws.Range("J116").Value = 2904
$ws.Range("K116").Value = 3475.125
$ws.Range("L116").Value = 8712
$ws.Range("M116").Value = -33.125
$ws.Range("N116").Value = -15596
$ws.Range("H117").Value = 5369.75
$ws.Range("I117").Value = 5264.5
$ws.Range("J117").Value = 5475
$ws.Range("K117").Value = 15793.5
$ws.Range("L117").Value = 16425
$ws.Range("M117").Value = -12351.5
$ws.Range("N117").Value = -23309
$ws.Range("H125").Value = 2710.9524
$ws.Range("I125").Value = 911.8182
$ws.Range("J125").Value = 4690
$ws.Range("K125").Value = 2735.4546
$ws.Range("L125").Value = 14070
$ws.Range("M125").Value = 2184.5454
$ws.Range("N125").Value = -23910

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 886.0769
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 968.7778
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 2906.3334
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -6746.3334
Write-Output "Shinryu_Profits sheets updated."
